$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K132").Value = 5643
$ws.Range("H132").Value = 1832.8572
$ws.Range("I132").Value = 1881
$ws.Range("M132").Value = -3113
$ws.Range("J132").Value = 1544
$ws.Range("N132").Value = -9692
$ws.Range("L132").Value = 4632
$ws.Range("N137").Value = -10275
$ws.Range("L137").Value = 5175
$ws.Range("J137").Value = 1725
$ws.Range("H137").Value = 1385.5714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2267287.5
$ws.Range("K32").Value = 2802774
$ws.Range("I32").Value = 2802774
$ws.Range("M32").Value = -2802487
$ws.Range("N61").Value = -4824
$ws.Range("H61").Value = 3270
$ws.Range("J61").Value = 4400
$ws.Range("K61").Value = 2987.5
$ws.Range("M61").Value = -2775.5
$ws.Range("L61").Value = 4400
$ws.Range("I61").Value = 2987.5
$ws.Range("I88").Value = 2532.4546
$ws.Range("L88").Value = 2710.8572
$ws.Range("K88").Value = 2532.4546
$ws.Range("N88").Value = -3522.8572
$ws.Range("H88").Value = 2601.8333
$ws.Range("J88").Value = 2710.8572
$ws.Range("M88").Value = -2126.4546
$ws.Range("J91").Value = 2710.8572
$ws.Range("H91").Value = 2601.8333
$ws.Range("K91").Value = 2532.4546
$ws.Range("N91").Value = -5518.8572
$ws.Range("L91").Value = 2710.8572
$ws.Range("M91").Value = -1128.4546
$ws.Range("I91").Value = 2532.4546
$ws.Range("H102").Value = 1644.8572
$ws.Range("M102").Value = -22.85719999999992
$ws.Range("I102").Value = 1644.8572
$ws.Range("K102").Value = 1644.8572
$ws.Range("K132").Value = 9641.000100000001
$ws.Range("H132").Value = 3213.6667
$ws.Range("I132").Value = 3213.6667
$ws.Range("M132").Value = -7111.000100000001
$ws.Range("J132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("L132").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("L136").Value = 13200
$ws.Range("I136").Value = 2987.5
$ws.Range("K136").Value = 8962.5
$ws.Range("H136").Value = 3270
$ws.Range("J136").Value = 4400
$ws.Range("M136").Value = -6412.5
$ws.Range("N136").Value = -18300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("J10").Value = 0
$ws.Range("N54").Value = -5968
$ws.Range("J54").Value = 5000
$ws.Range("L54").Value = 5000
$ws.Range("H54").Value = 5000
$ws.Range("N64").Value = -1298.4
$ws.Range("L64").Value = 848.4
$ws.Range("H64").Value = 876.8570999999999
$ws.Range("J64").Value = 848.4
$ws.Range("J67").Value = 848.4
$ws.Range("H67").Value = 876.8570999999999
$ws.Range("N67").Value = -2408.4
$ws.Range("L67").Value = 848.4
$ws.Range("H134").Value = 2769.5715
$ws.Range("M134").Value = -5773.7145
$ws.Range("K134").Value = 8308.7145
$ws.Range("I134").Value = 2769.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1659.6
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1224
$ws.Range("J12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("N12").Value = -5340
$ws.Range("H12").Value = 2502.5
$ws.Range("I31").Value = 1804
$ws.Range("H31").Value = 1804
$ws.Range("M31").Value = -1509
$ws.Range("K31").Value = 1804
$ws.Range("H34").Value = 1804
$ws.Range("M34").Value = -1602
$ws.Range("I34").Value = 1804
$ws.Range("K34").Value = 1804
$ws.Range("H134").Value = 2099.8
$ws.Range("M134").Value = -3465
$ws.Range("N134").Value = -12567
$ws.Range("J134").Value = 2499
$ws.Range("L134").Value = 7497
$ws.Range("K134").Value = 6000
$ws.Range("I134").Value = 2000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K12").Value = 789
$ws.Range("J12").Value = 230.4
$ws.Range("L12").Value = 691.2
$ws.Range("N12").Value = -1037.2
$ws.Range("M12").Value = -616
$ws.Range("I12").Value = 263
$ws.Range("H12").Value = 252.13333
$ws.Range("K80").Value = 11625
$ws.Range("I80").Value = 3875
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("M80").Value = -10689
$ws.Range("J80").Value = 0
$ws.Range("H80").Value = 3875
$ws.Range("M83").Value = -30195
$ws.Range("L83").Value = 0
$ws.Range("I83").Value = 3875
$ws.Range("H83").Value = 3875
$ws.Range("K83").Value = 34875
$ws.Range("J83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("L109").Value = 13200
$ws.Range("I109").Value = 530
$ws.Range("M109").Value = -550
$ws.Range("N109").Value = -15280
$ws.Range("J109").Value = 4400
$ws.Range("H109").Value = 1175
$ws.Range("K109").Value = 1590
$ws.Range("H115").Value = 3057.1428
$ws.Range("N115").Value = -14350
$ws.Range("J115").Value = 4000
$ws.Range("L115").Value = 12000
$ws.Range("K132").Value = 44995.5
$ws.Range("H132").Value = 5899.3335
$ws.Range("I132").Value = 4999.5
$ws.Range("M132").Value = -42465.5
$ws.Range("K141").Value = 28272
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -23092
$ws.Range("H141").Value = 8792
$ws.Range("J141").Value = 5000
$ws.Range("I141").Value = 9424
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L10").Value = 13500
$ws.Range("H10").Value = 16675167
$ws.Range("N10").Value = -13838
$ws.Range("J10").Value = 13500
$ws.Range("I31").Value = 1093
$ws.Range("H31").Value = 1093
$ws.Range("M31").Value = -801
$ws.Range("K31").Value = 1093
$ws.Range("K37").Value = 1093
$ws.Range("H37").Value = 1093
$ws.Range("I37").Value = 1093
$ws.Range("M37").Value = -816
$ws.Range("K132").Value = 7355.499899999999
$ws.Range("H132").Value = 2451.8333
$ws.Range("I132").Value = 2451.8333
$ws.Range("M132").Value = -4825.499899999999
$ws.Range("J140").Value = 142931.83
$ws.Range("L140").Value = 142931.83
$ws.Range("H140").Value = 142931.83
$ws.Range("N140").Value = -153291.83
$ws.Range("L141").Value = 59999
$ws.Range("H141").Value = 59999
$ws.Range("J141").Value = 59999
$ws.Range("N141").Value = -70359

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M22").Value = -804.75
$ws.Range("K22").Value = 1099.75
$ws.Range("H22").Value = 1233
$ws.Range("I22").Value = 1099.75
$ws.Range("I27").Value = 1099.75
$ws.Range("H27").Value = 1233
$ws.Range("M27").Value = -992.75
$ws.Range("K27").Value = 1099.75
$ws.Range("K46").Value = 1797.4
$ws.Range("N46").Value = -3993.2
$ws.Range("L46").Value = 3617.2
$ws.Range("I46").Value = 1797.4
$ws.Range("J46").Value = 3617.2
$ws.Range("M46").Value = -1609.4
$ws.Range("H46").Value = 2707.3
$ws.Range("K55").Value = 1499.875
$ws.Range("H55").Value = 1766.75
$ws.Range("I55").Value = 1499.875
$ws.Range("M55").Value = -1326.875
$ws.Range("K132").Value = 20275.6362
$ws.Range("H132").Value = 5872.0625
$ws.Range("I132").Value = 6758.5454
$ws.Range("M132").Value = -17745.6362
$ws.Range("I136").Value = 1604.625
$ws.Range("K136").Value = 4813.875
$ws.Range("H136").Value = 2283.6
$ws.Range("M136").Value = -2263.875
$ws.Range("H139").Value = 90650
$ws.Range("K139").Value = 90650
$ws.Range("M139").Value = -85510
$ws.Range("I139").Value = 90650

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K132").Value = 3629.4546
$ws.Range("H132").Value = 1680.5385
$ws.Range("I132").Value = 1209.8182
$ws.Range("M132").Value = -1099.4546
$ws.Range("J132").Value = 4269.5
$ws.Range("N132").Value = -17868.5
$ws.Range("L132").Value = 12808.5
$ws.Range("I136").Value = 2846.4
$ws.Range("K136").Value = 8539.200000000001
$ws.Range("H136").Value = 2807.7144
$ws.Range("M136").Value = -5989.200000000001
